# Auto-generated edit script applying the XLSX diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2952.6365
$ws.Range("I53").Value = 470
$ws.Range("J53").Value = 3504.3333
$ws.Range("K53").Value = 470
$ws.Range("L53").Value = 3504.3333
$ws.Range("M53").Value = 167
$ws.Range("N53").Value = -4778.3333

$ws.Range("H106").Value = 7094637.5
$ws.Range("I106").Value = 12347357
$ws.Range("K106").Value = 12347357
$ws.Range("M106").Value = -12346726

$ws.Range("H116").Value = 6491.75
$ws.Range("I116").Value = 3668.3333
$ws.Range("J116").Value = 7432.8887
$ws.Range("K116").Value = 3668.3333
$ws.Range("L116").Value = 7432.8887
$ws.Range("M116").Value = -226.3332999999998
$ws.Range("N116").Value = -14316.8887

$ws.Range("H135").Value = 15630049
$ws.Range("I135").Value = 640.5
$ws.Range("K135").Value = 5764.5
$ws.Range("M135").Value = -3229.5

$ws.Range("H137").Value = 1789.2667
$ws.Range("I137").Value = 1454.1111
$ws.Range("K137").Value = 4362.3333
$ws.Range("M137").Value = -1812.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4000
$ws.Range("I2").Value = 3500
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 3500
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -3387
$ws.Range("N2").Value = -5226

$ws.Range("H61").Value = 2905.32
$ws.Range("I61").Value = 3242.45
$ws.Range("J61").Value = 1556.8
$ws.Range("K61").Value = 3242.45
$ws.Range("L61").Value = 1556.8
$ws.Range("M61").Value = -3030.45
$ws.Range("N61").Value = -1980.8

$ws.Range("H63").Value = 4466251.5
$ws.Range("I63").Value = 2352
$ws.Range("J63").Value = 15626000
$ws.Range("K63").Value = 2352
$ws.Range("L63").Value = 15626000
$ws.Range("M63").Value = -1666
$ws.Range("N63").Value = -15627372

$ws.Range("H66").Value = 4466251.5
$ws.Range("I66").Value = 2352
$ws.Range("J66").Value = 15626000
$ws.Range("K66").Value = 11760
$ws.Range("L66").Value = 78130000
$ws.Range("M66").Value = -8328
$ws.Range("N66").Value = -78136864

$ws.Range("H74").Value = 43479884
$ws.Range("I74").Value = 71429170
$ws.Range("J74").Value = 3221.4443
$ws.Range("K74").Value = 71429170
$ws.Range("L74").Value = 3221.4443
$ws.Range("M74").Value = -71428296
$ws.Range("N74").Value = -4969.4443

$ws.Range("H77").Value = 43479884
$ws.Range("I77").Value = 71429170
$ws.Range("J77").Value = 3221.4443
$ws.Range("K77").Value = 357145850
$ws.Range("L77").Value = 16107.2215
$ws.Range("M77").Value = -357141482
$ws.Range("N77").Value = -24843.2215

$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 3500
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3500
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1206
$ws.Range("N116").Value = -9588

$ws.Range("H136").Value = 2905.32
$ws.Range("I136").Value = 3242.45
$ws.Range("J136").Value = 1556.8
$ws.Range("K136").Value = 9727.349999999999
$ws.Range("L136").Value = 4670.4
$ws.Range("M136").Value = -7177.349999999999
$ws.Range("N136").Value = -9770.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4000
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 3500
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -3386
$ws.Range("N3").Value = -5228

$ws.Range("H99").Value = 1565.2941
$ws.Range("I99").Value = 1246.6364
$ws.Range("J99").Value = 2149.5
$ws.Range("K99").Value = 1246.6364
$ws.Range("L99").Value = 2149.5
$ws.Range("M99").Value = 251.3635999999999
$ws.Range("N99").Value = -5145.5

$ws.Range("H134").Value = 2862.2683
$ws.Range("I134").Value = 3164.5
$ws.Range("J134").Value = 1394.2858
$ws.Range("K134").Value = 9493.5
$ws.Range("L134").Value = 4182.857400000001
$ws.Range("M134").Value = -6958.5
$ws.Range("N134").Value = -9252.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1441.3572
$ws.Range("I5").Value = 1030.8
$ws.Range("J5").Value = 1669.4445
$ws.Range("K5").Value = 3092.4
$ws.Range("L5").Value = 5008.333500000001
$ws.Range("M5").Value = -2980.4
$ws.Range("N5").Value = -5232.333500000001

$ws.Range("H22").Value = 25675
$ws.Range("I22").Value = 50550
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 151650
$ws.Range("L22").Value = 2400
$ws.Range("M22").Value = -151481
$ws.Range("N22").Value = -2738

$ws.Range("H27").Value = 25675
$ws.Range("I27").Value = 50550
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 151650
$ws.Range("L27").Value = 2400
$ws.Range("M27").Value = -151548
$ws.Range("N27").Value = -2604

$ws.Range("H122").Value = 627.8333
$ws.Range("I122").Value = 326
$ws.Range("J122").Value = 714.0714
$ws.Range("K122").Value = 2934
$ws.Range("L122").Value = 6426.6426
$ws.Range("M122").Value = -484
$ws.Range("N122").Value = -11326.6426

$ws.Range("H131").Value = 106045.945
$ws.Range("I131").Value = 645.5714
$ws.Range("J131").Value = 114430.07
$ws.Range("K131").Value = 1936.7142
$ws.Range("L131").Value = 343290.21
$ws.Range("M131").Value = 3103.2858
$ws.Range("N131").Value = -353370.21

$ws.Range("H135").Value = 1441.3572
$ws.Range("I135").Value = 1030.8
$ws.Range("J135").Value = 1669.4445
$ws.Range("K135").Value = 9277.199999999999
$ws.Range("L135").Value = 15025.0005
$ws.Range("M135").Value = -6742.199999999999
$ws.Range("N135").Value = -20095.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 33338000
$ws.Range("J52").Value = 33338000
$ws.Range("L52").Value = 33338000
$ws.Range("N52").Value = -33338518

$ws.Range("H80").Value = 4335.294
$ws.Range("I80").Value = 3725
$ws.Range("J80").Value = 4523.077
$ws.Range("K80").Value = 3725
$ws.Range("L80").Value = 4523.077
$ws.Range("M80").Value = -2727
$ws.Range("N80").Value = -6519.077

$ws.Range("H83").Value = 4335.294
$ws.Range("I83").Value = 3725
$ws.Range("J83").Value = 4523.077
$ws.Range("K83").Value = 18625
$ws.Range("L83").Value = 22615.385
$ws.Range("M83").Value = -13633
$ws.Range("N83").Value = -32599.385

$ws.Range("H97").Value = 1822
$ws.Range("I97").Value = 1761.5333
$ws.Range("J97").Value = 1922.7778
$ws.Range("K97").Value = 1761.5333
$ws.Range("L97").Value = 1922.7778
$ws.Range("M97").Value = -1265.5333
$ws.Range("N97").Value = -2914.7778

$ws.Range("H102").Value = 26318882
$ws.Range("I102").Value = 33336298
$ws.Range("K102").Value = 33336298
$ws.Range("M102").Value = -33334676

$ws.Range("H126").Value = 5428.6206
$ws.Range("J126").Value = 6586.923
$ws.Range("L126").Value = 19760.769
$ws.Range("N126").Value = -24700.769

$ws.Range("H132").Value = 34192.312
$ws.Range("I132").Value = 2842.125
$ws.Range("J132").Value = 65542.5
$ws.Range("K132").Value = 8526.375
$ws.Range("L132").Value = 196627.5
$ws.Range("M132").Value = -5996.375
$ws.Range("N132").Value = -201687.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4000
$ws.Range("I82").Value = 4000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 4000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3639
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 4000
$ws.Range("I85").Value = 4000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 4000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2752
$ws.Range("N85").ClearContents()

$ws.Range("H122").Value = 728996.75
$ws.Range("I122").Value = 1785057.5
$ws.Range("J122").Value = 2955
$ws.Range("K122").Value = 5355172.5
$ws.Range("L122").Value = 8865
$ws.Range("M122").Value = -5352722.5
$ws.Range("N122").Value = -13765

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 38465316
$ws.Range("I136").Value = 58825476
$ws.Range("J136").Value = 7236.6665
$ws.Range("K136").Value = 176476428
$ws.Range("L136").Value = 21709.9995
$ws.Range("M136").Value = -176473878
$ws.Range("N136").Value = -26809.9995
